$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Processes sheet: swap the "description" (B) and "type" (E) columns.
#    A temporary scratch column (H) is used so the exchange is a clean swap
#    of both values AND formatting (Cut preserves the source cell's style).
# ---------------------------------------------------------------------------
$wsProc = $wb.Worksheets.Item("Processes")

$wsProc.Range("B1:B12").Cut($wsProc.Range("H1:H12")) | Out-Null
$wsProc.Range("E1:E12").Cut($wsProc.Range("B1:B12")) | Out-Null
$wsProc.Range("H1:H12").Cut($wsProc.Range("E1:E12")) | Out-Null
$wsProc.Range("H1:H12").Clear() | Out-Null

# Resize columns B and E to match their new (swapped) content.
$wsProc.Columns("B:B").ColumnWidth = 12.9
$wsProc.Columns("E:E").ColumnWidth = 19.2

# ---------------------------------------------------------------------------
# 2. WasteDefinition sheet: add a "recycle" column (C) with default values.
# ---------------------------------------------------------------------------
$wsWaste = $wb.Worksheets.Item("WasteDefinition")

$wsWaste.Cells.Item(1, 1).Copy($wsWaste.Cells.Item(1, 3)) | Out-Null
$wsWaste.Cells.Item(1, 3).Value = "recycle"
$wsWaste.Cells.Item(2, 3).Value = 0
$wsWaste.Cells.Item(3, 3).Value = 0
$wsWaste.Range("C2").Select() | Out-Null

# ---------------------------------------------------------------------------
# 3. Fix the worksheet-scoped defined name so it matches the new Processes
#    layout (4 columns instead of 5).
# ---------------------------------------------------------------------------
$defName = $wsProc.Names.Item("cgam_processes")
$defName.RefersTo = "=Processes!`$A`$1:`$D`$1"

# ---------------------------------------------------------------------------
# 4. Make Processes the active sheet/selection (must be done last so the
#    tabSelected / activeTab bookkeeping ends up on the right sheet).
# ---------------------------------------------------------------------------
$wsProc.Range("B1:B12").Select() | Out-Null
